# Apply updated "想去人数" (interest count) and "最低票价" (lowest ticket price)
# values to the 杭州-漫展信息 workbook, matching the upstream data refresh.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 5923
$wsExhibit.Range("F3").Value = 561
$wsExhibit.Range("F4").Value = 1124
$wsExhibit.Range("F5").Value = 1067
$wsExhibit.Range("F6").Value = 856
$wsExhibit.Range("F7").Value = 90
$wsExhibit.Range("F8").Value = 53
$wsExhibit.Range("F9").Value = 620
$wsExhibit.Range("F10").Value = 67
$wsExhibit.Range("F13").Value = 2086
$wsExhibit.Range("F15").Value = 1155
$wsExhibit.Range("G16").Value = "已售罄"
$wsExhibit.Range("F17").Value = 214
$wsExhibit.Range("F18").Value = 450
$wsExhibit.Range("F19").Value = 674
$wsExhibit.Range("F20").Value = 239
$wsExhibit.Range("F24").Value = 3800
$wsExhibit.Range("F27").Value = 111
$wsExhibit.Range("F30").Value = 550
$wsExhibit.Range("F32").Value = 56
$wsExhibit.Range("F33").Value = 26
$wsExhibit.Range("F35").Value = 333
$wsExhibit.Range("F36").Value = 866
$wsExhibit.Range("F38").Value = 73
$wsExhibit.Range("F39").Value = 91
$wsExhibit.Range("F40").Value = 95

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 757
$wsShow.Range("F6").Value = 414

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5923
$wsAll.Range("F4").Value = 561
$wsAll.Range("F5").Value = 1124
$wsAll.Range("F7").Value = 757
$wsAll.Range("F8").Value = 1067
$wsAll.Range("F9").Value = 856
$wsAll.Range("F11").Value = 414
$wsAll.Range("F12").Value = 90
$wsAll.Range("F13").Value = 53
$wsAll.Range("F14").Value = 620
$wsAll.Range("F15").Value = 67
$wsAll.Range("F19").Value = 2086
$wsAll.Range("F21").Value = 1155
$wsAll.Range("G22").Value = "已售罄"
$wsAll.Range("F23").Value = 214
$wsAll.Range("F24").Value = 450
$wsAll.Range("F26").Value = 674
$wsAll.Range("F27").Value = 239
$wsAll.Range("F30").Value = 3800
$wsAll.Range("F33").Value = 111
$wsAll.Range("F36").Value = 550
$wsAll.Range("F38").Value = 56
$wsAll.Range("F39").Value = 26
$wsAll.Range("F41").Value = 333
$wsAll.Range("F42").Value = 866
$wsAll.Range("F44").Value = 73
$wsAll.Range("F45").Value = 91
$wsAll.Range("F46").Value = 95

